$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift the last three header labels over by one column
# (C1/D1/E1 = max/prediction/rejection-f -> prediction/rejection-f/max)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data row: C2 becomes the text that used to be in D2/E2, E2 becomes numeric 1
$ws.Range("C2").Value = "f__Treponemataceae"
$ws.Range("D2").Value = "f__Treponemataceae"
$ws.Range("E2").Value = 1
